$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F column (想去人数) for several rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 318
$wsExhibit.Range("F20").Value = 669
$wsExhibit.Range("F26").Value = 2384
$wsExhibit.Range("F27").Value = 4926
$wsExhibit.Range("F31").Value = 1264
$wsExhibit.Range("F36").Value = 75

# Sheet "全部类型" - update F column (想去人数) for the corresponding rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 318
$wsAll.Range("F21").Value = 669
$wsAll.Range("F27").Value = 2384
$wsAll.Range("F28").Value = 4926
$wsAll.Range("F32").Value = 1264
$wsAll.Range("F37").Value = 75
